$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New province order + updated figures for rows 2-13 (A:E)
$data = @(
    @("Hanoi",      332400, 101600, 1175500, 5.8),
    @("Vinh Phuc",  123800,  30800,  331200, 5.65),
    @("Bac Ninh",    82300,  36300,  439400, 6.04),
    @("Quang Ninh", 610200,  17200,  211300,  4.9),
    @("Hai Duong",  165600,  63000,  742600, 5.94),
    @("Hai Phong",  152700,  37500,  484700, 6.29),
    @("Hung Yen",    92600,  39500,  489600, 6.21),
    @("Thai Binh",  157100,  80500, 1061900, 6.56),
    @("Ha Nam",      86200,  33400,  401600, 5.96),
    @("Nam Dinh",   165300,  76300,  937700, 6.05),
    @("Ninh Binh",  137800,  41800,  484300, 6.02),
    @("Bac Giang",  385000,  53800,  626600, 5.55)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]

    # Rows were trimmed slightly (15 -> 13.8) in the refreshed table
    $ws.Rows.Item($row).RowHeight = 13.8
}

# Columns were narrowed a touch when the table was refreshed
# (ColumnWidth is entered net of Excel's ~5/6-character cell-padding, which
# gets added back in when the grid width is written out on save)
$ws.Columns.Item(1).ColumnWidth = 12.258503401360565
$ws.Columns.Item(2).ColumnWidth = 12.396258503401366
$ws.Columns.Item(3).ColumnWidth = 9.559523809523766
$ws.Columns.Item(4).ColumnWidth = 11.314625850340166
$ws.Columns.Item(5).ColumnWidth = 15.906462585033966

# Touch the sheet's true last row so the grid extent grows along with the refresh
$ws.Rows.Item(1048576).RowHeight = 12.8

$ws.Range("C18").Select()
